$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Rename the TOC bookmark on the "puzzles-cloud" Heading 1
#    paragraph: _Toc16364600457675906173036202 -> _Toc16364603065534479712892223
# ------------------------------------------------------------------
$bm1 = $d.Bookmarks.Item("_Toc16364600457675906173036202")
$bm1Range = $bm1.Range
$bm1.Delete()
$d.Bookmarks.Add("_Toc16364603065534479712892223", $bm1Range)

# ------------------------------------------------------------------
# 2) Rename the TOC bookmark on the "mnogo je dobro bilo" Heading 2
#    paragraph: _Toc16364600458075441435261037 -> _Toc16364603065854199637577553
# ------------------------------------------------------------------
$bm2 = $d.Bookmarks.Item("_Toc16364600458075441435261037")
$bm2Range = $bm2.Range
$bm2.Delete()
$d.Bookmarks.Add("_Toc16364603065854199637577553", $bm2Range)

# ------------------------------------------------------------------
# 3) Insert a brand-new Heading 3 paragraph, "asd", right after the
#    "mnogo je dobro bilo" Heading 2 paragraph - mirroring the same
#    shape as the existing headings (heading text followed by two
#    blank lines, with a new TOC bookmark wrapping the whole thing).
# ------------------------------------------------------------------
$headingPara = $d.Bookmarks.Item("_Toc16364603065854199637577553").Range.Paragraphs(1)
$newRange = $headingPara.Range.InsertParagraphAfter()

$newPara = $headingPara.Next()
$newPara.Style = "Heading 3"
$newPara.Range.Text = "asd"

# Append two more blank lines as literal newline characters inside the
# same paragraph (matching the "`n`n" tail used by the other headings).
$newPara = $headingPara.Next()
$tailPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$tailPoint.Text = "`n"

$newPara = $headingPara.Next()
$tailPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$tailPoint.Text = "`n"

# Bookmark the whole new paragraph's text (minus the paragraph mark),
# exactly like the other heading bookmarks.
$newPara = $headingPara.Next()
$bm3Range = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$d.Bookmarks.Add("_Toc16364603066388608835155560", $bm3Range)

Write-Output "Edit applied."
